$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Activate()

# --- Row 60/61: HTS.Index.Pos.T / HTS.Index.Neg.T ---
# New disagg type + uid (newly added DEGS "ActiveIndex")
$ws.Range("G60:G61").Value2 = "ActiveIndex/Age/Sex/Result"
$ws.Range("H60:H61").Value2 = "blSWpCnfpCK"
$ws.Range("M60").Value2 = "tss_NewPos"
$ws.Range("M61").Value2 = "tss_NewNeg"

# --- Row 64/65: HTS_TST.ActiveOther.Pos.T / HTS_TST.ActiveOther.Neg.T ---
# New disagg type + uid (newly added DEGS "ActiveOther")
$ws.Range("G64:G65").Value2 = "ActiveOther/Age/Sex/Result"
$ws.Range("H64:H65").Value2 = "WcpRsApp6sL"
$ws.Range("M64").Value2 = "tss_Pos"
$ws.Range("M65").Value2 = "tss_Neg"

# --- Row 67/68: HTS_TST.Other.Pos.T / HTS_TST.Other.Neg.T ---
# New disagg type + uid (newly added DEGS "OtherFacility") replacing the old
# multi-modality concatenated disagg string
$ws.Range("G67:G68").Value2 = "OtherFacility/Age/Sex/Result"
$ws.Range("H67:H68").Value2 = "QCZnSIBW2BI"

# Restore the view's active cell/selection to match the authored state
$ws.Range("H73").Select()
